# Daily attendance processing - 2025-11-23 11:45:37
# Normalises the "Recorded By" (column G) cell values: the list of
# recorders in each cell is stored in reverse-chronological order, so we
# flip the comma-separated entries end-for-end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By"
$col = 7

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $raw = $cell.Value()

    if ($raw -eq $null) { continue }
    if ($raw -eq "") { continue }

    $parts = $raw -split ',\s*'
    $n = $parts.Length

    if ($n -gt 1) {
        $reversed = $parts[($n - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        if ($newVal -ne $raw) {
            $cell.Value = $newVal
        }
    }
}
